# Add 2022-Q1 sheet with fund holdings data, insert it before the "总计" (total)
# summary sheet, and prepend a 2022-Q1 row to the "总计" summary table.

$wb = $excel.ActiveWorkbook

$zongjiBefore = $wb.Worksheets.Item("总计")
$template = $wb.Worksheets.Item("2021-Q4")

# --- 1. Create the new "2022-Q1" worksheet right before "总计" -------------
$newSheet = $wb.Worksheets.Add($zongjiBefore)
$newSheet.Name = "2022-Q1"

# NOTE: worksheet references in this host resolve by position, so after the
# insert the old "$zongjiBefore" handle now points at the new sheet instead
# of the total sheet. Re-fetch "总计" by name once the sheet list is final.
$zongji = $wb.Worksheets.Item("总计")

# Copy header row + column-A number formatting from an existing quarter sheet
# so the new sheet's styling (bold header, thin border, centered) matches.
$template.Range("A1:H1").Copy()
$newSheet.Range("A1:H1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$newSheet.Range("A2:A19").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Columns B-G hold text that looks numeric (fund codes with leading zeros,
# fixed-decimal percentages, ...); force text formatting up front so COM
# doesn't silently coerce the assigned strings into numbers.
$newSheet.Range("B2:G19").NumberFormat = "@"

$rows = @(
    @(0, "010902", "博时成长领航灵活配置混合A", "62.64", "73.89", "3.37", "2.1110", 9),
    @(1, "012366", "上投摩根安荣回报混合型证券投资基金A", "41.47", "21.90", "2.89", "1.1985", 1),
    @(2, "012463", "博时成长优势混合型证券投资基金A", "19.57", "75.80", "3.78", "0.7397", 7),
    @(3, "012367", "上投摩根安荣回报混合型证券投资基金C", "23.94", "21.90", "2.89", "0.6919", 1),
    @(4, "011033", "南方宝恒混合型证券投资基金A", "47.68", "20.19", "0.86", "0.4100", 4),
    @(5, "004823", "上投摩根安裕回报混合A", "11.29", "30.12", "3.46", "0.3906", 1),
    @(6, "011740", "博时成长精选混合A", "7.32", "75.84", "4.35", "0.3184", 6),
    @(7, "008966", "博时成长优选两年封闭运作灵活配置混合A", "7.62", "80.79", "4.04", "0.3078", 9),
    @(8, "010903", "博时成长领航灵活配置混合C", "8.01", "73.89", "3.37", "0.2699", 9),
    @(9, "004824", "上投摩根安裕回报混合C", "7.46", "30.12", "3.46", "0.2581", 1),
    @(10, "010742", "南方宁悦一年持有期混合A", "21.20", "23.31", "1.21", "0.2565", 1),
    @(11, "011034", "南方宝恒混合型证券投资基金C", "22.13", "20.19", "0.86", "0.1903", 4),
    @(12, "011741", "博时成长精选混合C", "1.39", "75.84", "4.35", "0.0605", 6),
    @(13, "012464", "博时成长优势混合型证券投资基金C", "0.83", "75.80", "3.78", "0.0314", 7),
    @(14, "010743", "南方宁悦一年持有期混合C", "2.16", "23.31", "1.21", "0.0261", 1),
    @(15, "008967", "博时成长优选两年封闭运作灵活配置混合C", "0.63", "80.79", "4.04", "0.0255", 9),
    @(16, "004316", "前海开源沪港深裕鑫灵活配置混合A", "0.64", "90.55", "3.06", "0.0196", 8),
    @(17, "004317", "前海开源沪港深裕鑫灵活配置混合C", "0.47", "90.55", "3.06", "0.0144", 8)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# --- 2. Prepend a 2022-Q1 row to the "总计" summary sheet -------------------
$zongji.Rows.Item(2).Insert()
$zongji.Range("A2:D2").ClearFormats()
# Row 3 (the old row 2, "2021-Q4") still carries the original "A column" style
# (bold/centered/bordered) - copy just that formatting onto the new A2.
$zongji.Range("A3").Copy()
$zongji.Range("A2").PasteSpecial(-4122)
$zongji.Range("A2").Value = 0
$zongji.Range("B2").Value = "2022-Q1"
$zongji.Range("C2").Value = 18
$zongji.Range("D2").Value = 7.32

# The A column is a running 0-based row index; renumber the rows pushed down
# by the insert (they kept their original 0..3 index values).
$zongji.Range("A3").Value = 1
$zongji.Range("A4").Value = 2
$zongji.Range("A5").Value = 3
$zongji.Range("A6").Value = 4
